# Updated cryptos list — refresh Price (D) and Volume(1h) (E) columns for
# rows 2-51 on the active sheet to match the latest scrape.
#
# Note: a few Price values round-trip cleanly as General numbers (Excel
# would store "230.57" fine), but some ("6.20", "2.60", "1.00", "552.50")
# would lose their trailing zero if Excel auto-typed them as numbers, so
# those specific cells are explicitly formatted as Text first (exactly
# what a user would do in Excel to keep the literal text) before the
# value is written, matching the original inline-string data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '91.739.24'
$ws.Range("E2").Value = '  -2.45%  '

$ws.Range("D3").Value = '3.319.90'
$ws.Range("E3").Value = '  -3.44%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").Value = '230.57'
$ws.Range("E5").Value = '  -2.81%  '

$ws.Range("D6").Value = '614.64'
$ws.Range("E6").Value = '  -3.54%  '

$ws.Range("E7").Value = '  -1.72%  '

$ws.Range("E8").Value = '  -2.94%  '

$ws.Range("E9").Value = '  +0.01%  '

$ws.Range("D10").Value = '0.947'
$ws.Range("E10").Value = '  -2.16%  '

$ws.Range("D11").Value = '3.318.42'
$ws.Range("E11").Value = '  -3.53%  '

$ws.Range("D12").Value = '42.21'
$ws.Range("E12").Value = '  +0.28%  '

$ws.Range("E13").Value = '  -1.72%  '

$ws.Range("D14").Value = '5.99'
$ws.Range("E14").Value = '  -2.42%  '

$ws.Range("D15").Value = '91.584.39'
$ws.Range("E15").Value = '  -2.38%  '

$ws.Range("D16").Value = '3.945.20'
$ws.Range("E16").Value = '  -3.28%  '

$ws.Range("E17").Value = '  -3.66%  '

$ws.Range("E18").Value = '  -3.74%  '

$ws.Range("D19").Value = '3.317.03'
$ws.Range("E19").Value = '  -3.54%  '

$ws.Range("D20").Value = '17.22'
$ws.Range("E20").Value = '  -1.99%  '

$ws.Range("D21").Value = '10.86'
$ws.Range("E21").Value = '  -3.78%  '

$ws.Range("E22").Value = '  +9.50%  '

$ws.Range("D23").Value = '490.71'
$ws.Range("E23").Value = '  -1.08%  '

$ws.Range("D24").Value = '0.449'
$ws.Range("E24").Value = '  -10.01%  '

$ws.Range("E25").Value = '  -3.32%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.20'
$ws.Range("E26").Value = '  -5.49%  '

$ws.Range("D27").Value = '91.37'
$ws.Range("E27").Value = '  +0.61%  '

$ws.Range("E28").Value = '  -1.00%  '

$ws.Range("D29").Value = '3.497.99'
$ws.Range("E29").Value = '  -3.44%  '

$ws.Range("E30").Value = '  -0.24%  '

$ws.Range("E31").Value = '  -5.89%  '

$ws.Range("E32").Value = '  +2.39%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.60'
$ws.Range("E33").Value = '  -4.73%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.00'
$ws.Range("E34").Value = '  +0.29%  '

$ws.Range("D35").Value = '0.171'
$ws.Range("E35").Value = '  -5.55%  '

$ws.Range("E36").Value = '  -7.04%  '

$ws.Range("D37").Value = '0.523'
$ws.Range("E37").Value = '  -6.57%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '552.50'
$ws.Range("E38").Value = '  +2.34%  '

$ws.Range("E39").Value = '  -0.02%  '

$ws.Range("E40").Value = '  -4.02%  '

$ws.Range("E41").Value = '  -2.06%  '

$ws.Range("E42").Value = '  -5.62%  '

$ws.Range("E43").Value = '  -6.91%  '

$ws.Range("D44").Value = '23.67'
$ws.Range("E44").Value = '  -1.49%  '

$ws.Range("D45").Value = '3.59'
$ws.Range("E45").Value = '  +2.92%  '

$ws.Range("E46").Value = '  -0.97%  '

$ws.Range("E47").Value = '  -0.19%  '

$ws.Range("E48").Value = '  -2.38%  '

$ws.Range("D49").Value = '2.09'
$ws.Range("E49").Value = '  -2.24%  '

$ws.Range("D50").Value = '51.66'
$ws.Range("E50").Value = '  -3.46%  '

$ws.Range("E51").Value = '  -0.97%  '
